# Fruta / hortaliza, semanal
# A new weekly price-report row is inserted at row 141 (pushing the
# existing rows 141-213 down to 142-214), and the new row is populated
# with its own data set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 141, shifting rows 141:213 -> 142:214
$ws.Rows("141:141").Insert()

# Populate the newly inserted row 141 with the new record's data
$ws.Cells.Item(141, 1).Value  = 4
$ws.Cells.Item(141, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(141, 3).Value  = "Los Lagos"
$ws.Cells.Item(141, 4).Value  = 44572
$ws.Cells.Item(141, 5).Value  = 10
$ws.Cells.Item(141, 6).Value  = 100112040
$ws.Cells.Item(141, 7).Value  = "Cilantro"
$ws.Cells.Item(141, 8).Value  = "Sin especificar"
$ws.Cells.Item(141, 9).Value  = "Primera"
$ws.Cells.Item(141, 10).Value = 120
$ws.Cells.Item(141, 11).Value = 10000
$ws.Cells.Item(141, 12).Value = 10000
$ws.Cells.Item(141, 13).Value = 10000
$ws.Cells.Item(141, 14).Value = "`$/docena de atados (2 kilos)"
$ws.Cells.Item(141, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(141, 16).Value = 5000
$ws.Cells.Item(141, 17).Value = 2
$ws.Cells.Item(141, 18).Value = "Hortaliza"
